$d = $word.ActiveDocument

$pairs = @(
    @("72-68=", "95-32="),
    @("6+45=", "43+40="),
    @("64-31=", "39+44="),
    @("31+64=", "88-6="),
    @("77-75=", "78-14="),
    @("46+3=", "4-3="),
    @("27+66=", "42+26="),
    @("40+10=", "6+56="),
    @("11+57=", "38-36="),
    @("50+44=", "75-49="),
    @("15+2=", "65-39="),
    @("60+36=", "92-37="),
    @("42+9=", "57+30="),
    @("34+39=", "95-22="),
    @("99-14=", "88+4="),
    @("74-17=", "94-51="),
    @("70-60=", "26-16="),
    @("56-55=", "31+57="),
    @("84-43=", "29+64="),
    @("0+14=", "37+39="),
    @("35+59=", "7+29="),
    @("71+9=", "54-15="),
    @("63-23=", "24-22="),
    @("79-4=", "39+31="),
    @("46-38=", "77-64="),
    @("75-6=", "94+5="),
    @("86-73=", "21+11="),
    @("13+4=", "5+15="),
    @("32+47=", "12+30="),
    @("15+26=", "0+63="),
    @("47+13=", "12+22="),
    @("83+0=", "13+63="),
    @("19+39=", "17+9="),
    @("94-78=", "70+29="),
    @("79-54=", "39-31="),
    @("40-11=", "86-33="),
    @("59-6=", "81-13="),
    @("15-0=", "36+47="),
    @("53+13=", "76+14="),
    @("24+68=", "10+71="),
    @("23+7=", "91-2="),
    @("20+47=", "74-67="),
    @("86-53=", "38+24="),
    @("37+41=", "16+60="),
    @("67-65=", "86-26="),
    @("62-16=", "38-27="),
    @("73-57=", "54+31="),
    @("93-49=", "32+58="),
    @("99-60=", "50+27="),
    @("42-32=", "54-13="),
    @("52-32=", "52-4="),
    @("67+23=", "34-2="),
    @("71-23=", "9+15="),
    @("10+22=", "85-65="),
    @("40+44=", "76-9="),
    @("73-72=", "28+59="),
    @("73-31=", "43+0="),
    @("52-47=", "1+0="),
    @("5+52=", "19+11="),
    @("39+2=", "18-3="),
    @("87+0=", "30-21="),
    @("0+17=", "90-60="),
    @("46-32=", "54-37="),
    @("6+29=", "94-12="),
    @("90-56=", "56+3="),
    @("82-25=", "98-2="),
    @("81-40=", "23+76="),
    @("91+8=", "14+19="),
    @("71+10=", "17+55="),
    @("57-49=", "37+54="),
    @("99-53=", "11+25="),
    @("65-27=", "86-51="),
    @("1+43=", "7+31="),
    @("11+37=", "33-5="),
    @("37+12=", "70-18="),
    @("38-35=", "83-78="),
    @("85-67=", "81-80="),
    @("21+36=", "39+51="),
    @("78-46=", "97-8="),
    @("26+6=", "62+2="),
    @("48+15=", "0+77="),
    @("36+25=", "58+38="),
    @("43-37=", "10+40="),
    @("52+26=", "33+56="),
    @("2+66=", "92-74="),
    @("97-14=", "80-66="),
    @("22+3=", "4+55="),
    @("21+42=", "39+42="),
    @("75-51=", "42+48="),
    @("43+23=", "22+9="),
    @("74+10=", "98-6="),
    @("31-2=", "93-69="),
    @("14+33=", "52+47="),
    @("12+49=", "49+17="),
    @("18+1=", "34+65="),
    @("8+45=", "73-71="),
    @("6+84=", "89-26="),
    @("3+60=", "62-30="),
    @("72+22=", "8+29="),
    @("68+26=", "89-2="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
